$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "param_TimeStep_starting_index"
$ws.Cells.Item(2, 2).Value = 5
$ws.Cells.Item(3, 1).Value = "param_demand1_inv_cost_starting_index"
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(4, 1).Value = "param_demand1_op_cost_starting_index"
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(5, 1).Value = "param_demand2_op_cost_starting_index"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(6, 1).Value = "param_demand2_inv_cost_starting_index"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(7, 1).Value = "param_P_net1_demand2_starting_index"
$ws.Cells.Item(7, 2).Value = 500
$ws.Cells.Item(8, 1).Value = "param_Q_from_net1_starting_index"
$ws.Cells.Item(8, 2).Value = 1143.156692505584
$ws.Cells.Item(9, 1).Value = "param_net1_inv_cost_starting_index"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(10, 1).Value = "param_Q_net1_demand2_starting_index"
$ws.Cells.Item(10, 2).Value = 1000
$ws.Cells.Item(11, 1).Value = "param_P_to_net1_starting_index"
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(12, 1).Value = "param_P_net1_demand1_starting_index"
$ws.Cells.Item(12, 2).Value = 95.4377950037229
$ws.Cells.Item(13, 1).Value = "param_net1_emissions_starting_index"
$ws.Cells.Item(13, 2).Value = 649.0650187324193
$ws.Cells.Item(14, 1).Value = "param_P_net1_bat1_starting_index"
$ws.Cells.Item(14, 2).Value = 73.68258451606111
$ws.Cells.Item(15, 1).Value = "param_net1_buy_electric_starting_index"
$ws.Cells.Item(15, 2).Value = 267.6481518079136
$ws.Cells.Item(16, 1).Value = "param_net1_buy_thermal_starting_index"
$ws.Cells.Item(16, 2).Value = 400.1048423769545
$ws.Cells.Item(17, 1).Value = "param_P_from_net1_starting_index"
$ws.Cells.Item(17, 2).Value = 669.1203795197839
$ws.Cells.Item(18, 1).Value = "param_net1_sell_electric_starting_index"
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(19, 1).Value = "param_Q_to_net1_starting_index"
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(20, 1).Value = "param_net1_sell_thermal_starting_index"
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(21, 1).Value = "param_Q_net1_demand1_starting_index"
$ws.Cells.Item(21, 2).Value = 143.1566925055843
$ws.Cells.Item(22, 1).Value = "param_P_pv1_net1_starting_index"
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(23, 1).Value = "param_pv1_op_cost_starting_index"
$ws.Cells.Item(23, 2).Value = 1
$ws.Cells.Item(24, 1).Value = "param_P_pv1_bat1_starting_index"
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(25, 1).Value = "param_pv1_inv_cost_starting_index"
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(26, 1).Value = "param_P_pv1_demand1_starting_index"
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(27, 1).Value = "param_pv1_emissions_starting_index"
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(28, 1).Value = "param_P_from_pv1_starting_index"
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(29, 1).Value = "param_P_pv1_demand2_starting_index"
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(30, 1).Value = "param_bat1_emissions_starting_index"
$ws.Cells.Item(30, 2).Value = 0.2046738458779475
$ws.Cells.Item(31, 1).Value = "param_bat1_inv_cost_starting_index"
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(32, 1).Value = "param_bat1_SOC_max_starting_index"
$ws.Cells.Item(32, 2).Value = 0.9999845529025807
$ws.Cells.Item(33, 1).Value = "param_P_to_bat1_starting_index"
$ws.Cells.Item(33, 2).Value = 73.68258451606111
$ws.Cells.Item(34, 1).Value = "param_P_bat1_demand2_starting_index"
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(35, 1).Value = "param_bat1_K_ch_starting_index"
$ws.Cells.Item(35, 2).Value = 1
$ws.Cells.Item(36, 1).Value = "param_bat1_integer_starting_index"
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(37, 1).Value = "param_bat1_op_cost_starting_index"
$ws.Cells.Item(37, 2).Value = 1
$ws.Cells.Item(38, 1).Value = "param_bat1_K_dis_starting_index"
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(39, 1).Value = "param_P_bat1_demand1_starting_index"
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(40, 1).Value = "param_P_bat1_net1_starting_index"
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(41, 1).Value = "param_P_from_bat1_starting_index"
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(42, 1).Value = "param_bat1_SOC_starting_index"
$ws.Cells.Item(42, 2).Value = 0.9999845529025807
$ws.Cells.Item(43, 1).Value = "param_bat1_cumulated_aging_starting_index"
$ws.Cells.Item(43, 2).Value = 0.00001544709741934352
$ws.Cells.Item(44, 1).Value = "param_total_emissions_starting_index"
$ws.Cells.Item(44, 2).Value = 649.2696925782972
$ws.Cells.Item(45, 1).Value = "param_total_operation_cost_starting_index"
$ws.Cells.Item(45, 2).Value = 2
$ws.Cells.Item(46, 1).Value = "param_total_buy_starting_index"
$ws.Cells.Item(46, 2).Value = 667.752994184868
$ws.Cells.Item(47, 1).Value = "param_total_sell_starting_index"
$ws.Cells.Item(47, 2).Value = 0
